# Validacion de carga de datos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The country value stored in B1 was "Belize"; replace it with the ISO3 code "BLZ"
$ws.Range("B1").Value = "BLZ"

# Leave the selection on B1, matching the state after the edit was made
$ws.Range("B1").Select()
